# Update the math expressions in the worksheet table.
# Each cell is addressed directly via Table.Cell(row, column) to avoid
# ambiguity from duplicate expressions (e.g. "4+11=" appears twice, only
# one of which changes).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "14+28="
$t.Cell(1, 2).Range.Text = "53-33="
$t.Cell(1, 3).Range.Text = "5+40="
$t.Cell(1, 4).Range.Text = "28+4="
$t.Cell(1, 5).Range.Text = "54+37="
$t.Cell(2, 1).Range.Text = "37-22="
$t.Cell(2, 2).Range.Text = "34-5="
$t.Cell(2, 3).Range.Text = "3+95="
$t.Cell(2, 4).Range.Text = "40-19="
$t.Cell(2, 5).Range.Text = "34+15="
$t.Cell(3, 1).Range.Text = "35+18="
$t.Cell(3, 2).Range.Text = "82+15="
$t.Cell(3, 3).Range.Text = "27+56="
$t.Cell(3, 4).Range.Text = "24+71="
$t.Cell(3, 5).Range.Text = "49-46="
$t.Cell(4, 1).Range.Text = "23-7="
$t.Cell(4, 2).Range.Text = "81-23="
$t.Cell(4, 3).Range.Text = "50+1="
$t.Cell(4, 4).Range.Text = "6+92="
$t.Cell(4, 5).Range.Text = "85-44="
$t.Cell(5, 1).Range.Text = "28+56="
$t.Cell(5, 2).Range.Text = "67-8="
$t.Cell(5, 3).Range.Text = "14+65="
$t.Cell(5, 4).Range.Text = "46-32="
$t.Cell(5, 5).Range.Text = "72-4="
$t.Cell(6, 1).Range.Text = "73+9="
$t.Cell(6, 2).Range.Text = "93-12="
$t.Cell(6, 3).Range.Text = "53-28="
$t.Cell(6, 4).Range.Text = "55-16="
$t.Cell(6, 5).Range.Text = "23-15="
$t.Cell(7, 1).Range.Text = "52-40="
$t.Cell(7, 2).Range.Text = "97-58="
$t.Cell(7, 3).Range.Text = "17+61="
$t.Cell(7, 4).Range.Text = "38+40="
$t.Cell(7, 5).Range.Text = "96-76="
$t.Cell(8, 1).Range.Text = "61-5="
$t.Cell(8, 2).Range.Text = "90-17="
$t.Cell(8, 3).Range.Text = "7+29="
$t.Cell(8, 4).Range.Text = "74-17="
$t.Cell(8, 5).Range.Text = "81-8="
$t.Cell(9, 1).Range.Text = "78-22="
$t.Cell(9, 2).Range.Text = "64-38="
$t.Cell(9, 3).Range.Text = "2+46="
$t.Cell(9, 4).Range.Text = "54-46="
$t.Cell(9, 5).Range.Text = "7+34="
$t.Cell(10, 1).Range.Text = "23-17="
$t.Cell(10, 2).Range.Text = "69+29="
$t.Cell(10, 3).Range.Text = "37-20="
$t.Cell(10, 4).Range.Text = "67-36="
$t.Cell(11, 1).Range.Text = "42+19="
$t.Cell(11, 2).Range.Text = "84-55="
$t.Cell(11, 3).Range.Text = "4+63="
$t.Cell(11, 4).Range.Text = "78-38="
$t.Cell(11, 5).Range.Text = "31+11="
$t.Cell(12, 1).Range.Text = "24-21="
$t.Cell(12, 2).Range.Text = "64-43="
$t.Cell(12, 3).Range.Text = "20-3="
$t.Cell(12, 4).Range.Text = "29+27="
$t.Cell(12, 5).Range.Text = "66-57="
$t.Cell(13, 1).Range.Text = "38+19="
$t.Cell(13, 2).Range.Text = "6+19="
$t.Cell(13, 3).Range.Text = "93-39="
$t.Cell(13, 4).Range.Text = "31+45="
$t.Cell(13, 5).Range.Text = "5+13="
$t.Cell(14, 1).Range.Text = "86-67="
$t.Cell(14, 2).Range.Text = "53-20="
$t.Cell(14, 3).Range.Text = "53-23="
$t.Cell(14, 4).Range.Text = "3+90="
$t.Cell(14, 5).Range.Text = "83+5="
$t.Cell(15, 1).Range.Text = "90-73="
$t.Cell(15, 2).Range.Text = "3+0="
$t.Cell(15, 3).Range.Text = "32-19="
$t.Cell(15, 4).Range.Text = "91-34="
$t.Cell(15, 5).Range.Text = "22-4="
$t.Cell(16, 1).Range.Text = "6+72="
$t.Cell(16, 2).Range.Text = "55-22="
$t.Cell(16, 3).Range.Text = "78+8="
$t.Cell(16, 4).Range.Text = "61-16="
$t.Cell(16, 5).Range.Text = "87-48="
$t.Cell(17, 1).Range.Text = "60+7="
$t.Cell(17, 2).Range.Text = "40+1="
$t.Cell(17, 3).Range.Text = "44-37="
$t.Cell(17, 4).Range.Text = "41-34="
$t.Cell(17, 5).Range.Text = "20+36="
$t.Cell(18, 1).Range.Text = "16+28="
$t.Cell(18, 2).Range.Text = "4+76="
$t.Cell(18, 3).Range.Text = "29-8="
$t.Cell(18, 4).Range.Text = "96-23="
$t.Cell(18, 5).Range.Text = "17+62="
$t.Cell(19, 1).Range.Text = "75-54="
$t.Cell(19, 2).Range.Text = "97-83="
$t.Cell(19, 3).Range.Text = "57-31="
$t.Cell(19, 4).Range.Text = "9+36="
$t.Cell(19, 5).Range.Text = "6-6="
$t.Cell(20, 1).Range.Text = "87-10="
$t.Cell(20, 2).Range.Text = "60+16="
$t.Cell(20, 3).Range.Text = "35+29="
$t.Cell(20, 4).Range.Text = "95-83="
$t.Cell(20, 5).Range.Text = "76-7="
